$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,10
$data[0,0] = 13.57063651909873
$data[0,1] = 6.883135427878539
$data[0,2] = 14.04898882569214
$data[0,3] = 48.34100886830125
$data[0,4] = 59.97392342419536
$data[0,5] = 22.10904475649088
$data[0,6] = 0
$data[0,7] = 11.52380739942606
$data[0,8] = 22.83994665585294
$data[0,9] = 9.598492697885041
$data[1,0] = 13.53292697867383
$data[1,1] = 6.849936529828728
$data[1,2] = 14.05080937376609
$data[1,3] = 48.48362694200738
$data[1,4] = 60.1672530875631
$data[1,5] = 22.19868091358136
$data[1,6] = 0
$data[1,7] = 11.5528234621601
$data[1,8] = 22.44576138341752
$data[1,9] = 9.612206686268324
$data[2,0] = 13.51255438158723
$data[2,1] = 6.829901216823169
$data[2,2] = 14.05401963081219
$data[2,3] = 48.58467146948104
$data[2,4] = 60.30612398497656
$data[2,5] = 22.25864897341566
$data[2,6] = 0
$data[2,7] = 11.57215446707737
$data[2,8] = 22.20396445273511
$data[2,9] = 9.621413889922827
$data[3,0] = 13.50495682638736
$data[3,1] = 6.821828827036046
$data[3,2] = 14.05585476252507
$data[3,3] = 48.62921963216266
$data[3,4] = 60.36775126273667
$data[3,5] = 22.2843226759032
$data[3,6] = 0
$data[3,7] = 11.58041309209699
$data[3,8] = 22.10561096013299
$data[3,9] = 9.625364174791171
$data[4,0] = 13.50373793978071
$data[4,1] = 6.82049410643333
$data[4,2] = 14.05619132905687
$data[4,3] = 48.63681994994006
$data[4,4] = 60.37828756389474
$data[4,5] = 22.28866033075512
$data[4,6] = 0
$data[4,7] = 11.58180745048015
$data[4,8] = 22.08929393776363
$data[4,9] = 9.626032104706777
$data[5,0] = 13.51244905982467
$data[5,1] = 6.829791970805734
$data[5,2] = 14.0540422456404
$data[5,3] = 48.5852586330023
$data[5,4] = 60.30693476743284
$data[5,5] = 22.25899021753088
$data[5,6] = 0
$data[5,7] = 11.57226430264375
$data[5,8] = 22.20263712956537
$data[5,9] = 9.621466361461314
$data[6,0] = 13.55706028276488
$data[6,1] = 6.871618700660706
$data[6,2] = 14.04918276910868
$data[6,3] = 48.38737649114059
$data[6,4] = 60.03637775225113
$data[6,5] = 22.13892561226485
$data[6,6] = 0
$data[6,7] = 11.53349776596814
$data[6,8] = 22.70405548101972
$data[6,9] = 9.60305825192899
$data[7,0] = 13.66635762834423
$data[7,1] = 6.956212735561886
$data[7,2] = 14.05621755630578
$data[7,3] = 48.10699244018205
$data[7,4] = 59.66729930535205
$data[7,5] = 21.94277223528512
$data[7,6] = 0
$data[7,7] = 11.46949266321208
$data[7,8] = 23.68394400974778
$data[7,9] = 9.57318340110306
$data[8,0] = 13.75957700279564
$data[8,1] = 7.01966494668428
$data[8,2] = 14.07142570477653
$data[8,3] = 47.96758328926384
$data[8,4] = 59.49649557265017
$data[8,5] = 21.82284086130081
$data[8,6] = 0
$data[8,7] = 11.42978668080459
$data[8,8] = 24.39496977430714
$data[8,9] = 9.555002554527951
$data[9,0] = 13.80469556325298
$data[9,1] = 7.048755349217985
$data[9,2] = 14.08050855463736
$data[9,3] = 47.91880848796682
$data[9,4] = 59.44095610367962
$data[9,5] = 21.7735810365708
$data[9,6] = 0
$data[9,7] = 11.41331120269168
$data[9,8] = 24.71513009310653
$data[9,9] = 9.547544386876043
$data[10,0] = 13.82216186083122
$data[10,1] = 7.05979855057573
$data[10,2] = 14.08425758110617
$data[10,3] = 47.9024581723891
$data[10,4] = 59.42313860275837
$data[10,5] = 21.75569349373323
$data[10,6] = 0
$data[10,7] = 11.40730048322451
$data[10,8] = 24.83578538495441
$data[10,9] = 9.54483654851165
$data[11,0] = 13.81838339394651
$data[11,1] = 7.057419067929763
$data[11,2] = 14.08343642315455
$data[11,3] = 47.90588502211626
$data[11,4] = 59.42683255842783
$data[11,5] = 21.75951175598141
$data[11,6] = 0
$data[11,7] = 11.40858485084398
$data[11,8] = 24.80982770167976
$data[11,9] = 9.545414559266133
$data[12,0] = 13.80612494739444
$data[12,1] = 7.049663362596531
$data[12,2] = 14.08081079595375
$data[12,3] = 47.917420780487
$data[12,4] = 59.43942568233619
$data[12,5] = 21.77209403961395
$data[12,6] = 0
$data[12,7] = 11.4128121240925
$data[12,8] = 24.72506875857264
$data[12,9] = 9.547319280765286
$data[13,0] = 13.79866561385343
$data[13,1] = 7.04491617052844
$data[13,2] = 14.07924278183885
$data[13,3] = 47.92476320286561
$data[13,4] = 59.44755868237949
$data[13,5] = 21.77990094577362
$data[13,6] = 0
$data[13,7] = 11.41543116708781
$data[13,8] = 24.67307244372466
$data[13,9] = 9.548501125151764
$data[14,0] = 13.7566822108881
$data[14,1] = 7.017767921771711
$data[14,2] = 14.07087553351624
$data[14,3] = 47.97106671284151
$data[14,4] = 59.50057352529049
$data[14,5] = 21.8261670750787
$data[14,6] = 0
$data[14,7] = 11.43089532444745
$data[14,8] = 24.37397087967489
$data[14,9] = 9.555506275421441
$data[15,0] = 13.7316151623892
$data[15,1] = 7.001167495014761
$data[15,2] = 14.06629573536254
$data[15,3] = 48.00323314287289
$data[15,4] = 59.53879238983139
$data[15,5] = 21.85590974272913
$data[15,6] = 0
$data[15,7] = 11.4407885330231
$data[15,8] = 24.18955979722945
$data[15,9] = 9.560011487390616
$data[16,0] = 13.71745317012459
$data[16,1] = 6.991640921041546
$data[16,2] = 14.06386541180393
$data[16,3] = 48.02311193055079
$data[16,4] = 59.56285916556485
$data[16,5] = 21.87351531342712
$data[16,6] = 0
$data[16,7] = 11.44662821944766
$data[16,8] = 24.08318808503491
$data[16,9] = 9.562679256968574
$data[17,0] = 13.7127023862655
$data[17,1] = 6.9884192449134
$data[17,2] = 14.06307760771192
$data[17,3] = 48.03007869181609
$data[17,4] = 59.57136486772217
$data[17,5] = 21.87956172219614
$data[17,6] = 0
$data[17,7] = 11.44863109267525
$data[17,8] = 24.04712377686146
$data[17,9] = 9.563595667840284
$data[18,0] = 13.73425717140388
$data[18,1] = 7.002932445388156
$data[18,2] = 14.06676217634136
$data[18,3] = 47.99966630057245
$data[18,4] = 59.53450800415646
$data[18,5] = 21.85269197409642
$data[18,6] = 0
$data[18,7] = 11.43971992455333
$data[18,8] = 24.2092228761961
$data[18,9] = 9.559523986378419
$data[19,0] = 13.80971528948004
$data[19,1] = 7.051940701649665
$data[19,2] = 14.08157361981617
$data[19,3] = 47.91397481584328
$data[19,4] = 59.43563933885351
$data[19,5] = 21.76837749240164
$data[19,6] = 0
$data[19,7] = 11.4115642786993
$data[19,8] = 24.74998116927765
$data[19,9] = 9.546756661908953
$data[20,0] = 13.86124691824871
$data[20,1] = 7.084127129070201
$data[20,2] = 14.09305711772098
$data[20,3] = 47.87033112158282
$data[20,4] = 59.3897683008039
$data[20,5] = 21.71774041218928
$data[20,6] = 0
$data[20,7] = 11.39449297333941
$data[20,8] = 25.09995591400083
$data[20,9] = 9.539090799217162
$data[21,0] = 13.83354399335233
$data[21,1] = 7.066935999400268
$data[21,2] = 14.08676376377073
$data[21,3] = 47.89248923448676
$data[21,4] = 59.41252682692516
$data[21,5] = 21.74435619837719
$data[21,6] = 0
$data[21,7] = 11.40348255797683
$data[21,8] = 24.91351754395052
$data[21,9] = 9.54312028454231
$data[22,0] = 13.73306194165199
$data[22,1] = 7.002134457273844
$data[22,2] = 14.06655066698926
$data[22,3] = 48.0012745541115
$data[22,4] = 59.53643845361059
$data[22,5] = 21.85414515107446
$data[22,6] = 0
$data[22,7] = 11.44020256913596
$data[22,8] = 24.20033428408997
$data[22,9] = 9.559744143645414
$data[23,0] = 13.6344909831952
$data[23,1] = 6.933079720708851
$data[23,2] = 14.05254648179361
$data[23,3] = 48.1712157705357
$data[23,4] = 59.74965228163482
$data[23,5] = 21.99160673902803
$data[23,6] = 0
$data[23,7] = 11.48552204045525
$data[23,8] = 23.41994033252323
$data[23,9] = 9.580601740683544

$ws.Range("C2:L25").Value = $data

